# Generate Report for Handback
#
# - Status for the 917425fb-... item ("Ready for handoff") becomes
#   "Handback transform failed" everywhere it is shown (Overview sheet,
#   and the Status column on the zh-cn / de-de detail sheets).
# - The zh-cn and de-de detail sheets get an "Error Detail" message on
#   that same row explaining the handback/handoff file name mismatch.
# - The "Error Detail" column (P) is widened to fit the new text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handback transform failed"

$zhError = "Handback file name: uf104uht.2rx is different with handoff file name: 917425fb-7ed2-4ddd-b175-a10317cb8a15.917f7ea59385eb31ea264a1a318302fa6c386951.zh-cn."
$deError = "Handback file name: uf104uht.2rx is different with handoff file name: 917425fb-7ed2-4ddd-b175-a10317cb8a15.917f7ea59385eb31ea264a1a318302fa6c386951.de-de."

# --- Overview sheet: update the zh-cn / de-de status columns (E3:F3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn detail sheet: Status + Error Detail for the 917425fb row ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("P3").Value = $zhError
# Widen the Error Detail column (16th / P) so the new text fits (raw
# stored column width of 40; ColumnWidth applies a +5/6 character offset
# on save, so back it out here to land exactly on 40).
$wsZh.Columns.Item(16).ColumnWidth = 40 - (5/6)

# --- de-de detail sheet: Status + Error Detail for the 917425fb row ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = 40 - (5/6)
